# Remove the record for account 004377713 (DANIELI, balance 28672.04)
# from the "Export" sheet. Removing the row shifts all subsequent rows
# up by one, matching the source data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetAccount = "004377713"

# Locate the row holding the account number so the edit is resilient to
# the exact row position; fall back to the known row if Find is unavailable.
$rowToDelete = 4
$found = $ws.Cells.Find($targetAccount)
if ($found -ne $null) {
    $rowToDelete = $found.Row
}

$ws.Range("A" + $rowToDelete).EntireRow.Delete()
